$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5774
$ws.Range("F7").Value = 9816
$ws.Range("F10").Value = 3914
$ws.Range("F11").Value = 79
$ws.Range("F14").Value = 213
$ws.Range("F20").Value = 633
$ws.Range("F21").Value = 3931
$ws.Range("F24").Value = 5408
$ws.Range("F25").Value = 442
$ws.Range("F26").Value = 2133
$ws.Range("F28").Value = 364
$ws.Range("F29").Value = 8061
$ws.Range("F30").Value = 34
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 2213
$ws.Range("F33").Value = 2223
$ws.Range("F35").Value = 1323
$ws.Range("F38").Value = 281
$ws.Range("F39").Value = 254
$ws.Range("F41").Value = 1188
$ws.Range("F43").Value = 181
$ws.Range("F44").Value = 1355
$ws.Range("F45").Value = 2134
$ws.Range("F46").Value = 139
$ws.Range("F47").Value = 233

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 151
$ws.Range("F11").Value = 128

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 592
$ws.Range("F3").Value = 772

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 772
$ws.Range("F6").Value = 5774
$ws.Range("F8").Value = 3914
$ws.Range("F9").Value = 79
$ws.Range("F16").Value = 151
$ws.Range("F18").Value = 633
$ws.Range("F19").Value = 3931
$ws.Range("F23").Value = 5408
$ws.Range("F24").Value = 442
$ws.Range("F25").Value = 2133
$ws.Range("F27").Value = 364
$ws.Range("F28").Value = 8061
$ws.Range("F29").Value = 34
$ws.Range("F30").Value = 2213
$ws.Range("F31").Value = 2223
$ws.Range("F33").Value = 1323
$ws.Range("F35").Value = 281
$ws.Range("F36").Value = 254
$ws.Range("F38").Value = 1188
$ws.Range("F40").Value = 181
$ws.Range("F42").Value = 1355
$ws.Range("F44").Value = 2134
$ws.Range("F45").Value = 139
$ws.Range("F46").Value = 233
